# Apply the data updates described in the commit:
#  - Normalize a handful of mis-ordered "Station" diagnosis labels so the
#    ligament names are listed in a consistent order.
#  - Add PreviousThiel / PreviousDiagnostic columns (BC, BD) to the main
#    data table, populated from the small standalone table that used to
#    live at A11:C19.
#  - Remove that standalone table (rows 11-19) entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel a few station results to a consistent ligament ordering.
$ws.Range("Q2").Value = "ACL, MCL"
$ws.Range("K3").Value = "ACL, MCL"
$ws.Range("R3").Value = "MCL, PCL"
$ws.Range("R4").Value = "ACL, LCL"
$ws.Range("N5").Value = "ACL, MCL"
$ws.Range("Q5").Value = "ACL, LCL"
$ws.Range("R5").Value = "ACL, MCL"
$ws.Range("R9").Value = "MCL, PCL"

# New headers for the appended columns.
$ws.Range("BC1").Value = "PreviousThiel"
$ws.Range("BD1").Value = "PreviousDiagnostic"

# Values from the old standalone table (A12:C19), keyed by Subject number
# and re-aligned onto the corresponding row of the main table (A2:A9).
$ws.Range("BC2").Value = 8
$ws.Range("BD2").Value = 5

$ws.Range("BC3").Value = 10
$ws.Range("BD3").Value = 5

$ws.Range("BC4").Value = 10
$ws.Range("BD4").Value = 1

$ws.Range("BC5").Value = 10
$ws.Range("BD5").Value = 7

$ws.Range("BC6").Value = 5
$ws.Range("BD6").Value = 3

$ws.Range("BC7").Value = 10
$ws.Range("BD7").Value = 4

$ws.Range("BC8").Value = 8
$ws.Range("BD8").Value = 1

$ws.Range("BC9").Value = 10
$ws.Range("BD9").Value = 2

# Remove the now-redundant standalone table (rows 11-19) completely so the
# used range collapses back down to row 9.
$ws.Range("A11:C19").EntireRow.Delete()

# Match the selection left behind in the saved file.
$ws.Range("R11").Select()
